{"js": "const replacements = [\n  [\"500\u00f73=\", \"718\u00f78=\"],\n  [\"940\u00f73=\", \"127\u00f75=\"],\n  [\"810\u00f79=\", \"630\u00f76=\"],\n  [\"656\u00f73=\", \"446\u00f79=\"],\n  [\"550\u00f72=\", \"801\u00f76=\"],\n  [\"701\u00f75=\", \"251\u00f76=\"],\n  [\"203\u00f78=\", \"481\u00f78=\"],\n  [\"299\u00f76=\", \"943\u00f77=\"],\n  [\"574\u00f76=\", \"899\u00f79=\"],\n  [\"861\u00f75=\", \"323\u00f76=\"],\n  [\"528\u00f76=\", \"592\u00f78=\"],\n  [\"394\u00f72=\", \"590\u00f77=\"],\n  [\"950\u00f79=\", \"377\u00f77=\"],\n  [\"167\u00f79=\", \"386\u00f73=\"],\n  [\"723\u00f76=\", \"967\u00f72=\"],\n  [\"268\u00f76=\", \"901\u00f77=\"],\n  [\"151\u00f74=\", \"900\u00f75=\"],\n  [\"374\u00f77=\", \"291\u00f76=\"],\n  [\"375\u00f79=\", \"530\u00f72=\"],\n  [\"902\u00f74=\", \"470\u00f72=\"],\n  [\"262\u00f77=\", \"783\u00f73=\"],\n  [\"281\u00f72=\", \"959\u00f76=\"],\n  [\"413\u00f75=\", \"308\u00f78=\"],\n  [\"867\u00f72=\", \"649\u00f72=\"],\n  [\"396\u00f76=\", \"281\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nreturn 'ok';", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"500\u00f73=\", \"718\u00f78=\"),\n    @(\"940\u00f73=\", \"127\u00f75=\"),\n    @(\"810\u00f79=\", \"630\u00f76=\"),\n    @(\"656\u00f73=\", \"446\u00f79=\"),\n    @(\"550\u00f72=\", \"801\u00f76=\"),\n    @(\"701\u00f75=\", \"251\u00f76=\"),\n    @(\"203\u00f78=\", \"481\u00f78=\"),\n    @(\"299\u00f76=\", \"943\u00f77=\"),\n    @(\"574\u00f76=\", \"899\u00f79=\"),\n    @(\"861\u00f75=\", \"323\u00f76=\"),\n    @(\"528\u00f76=\", \"592\u00f78=\"),\n    @(\"394\u00f72=\", \"590\u00f77=\"),\n    @(\"950\u00f79=\", \"377\u00f77=\"),\n    @(\"167\u00f79=\", \"386\u00f73=\"),\n    @(\"723\u00f76=\", \"967\u00f72=\"),\n    @(\"268\u00f76=\", \"901\u00f77=\"),\n    @(\"151\u00f74=\", \"900\u00f75=\"),\n    @(\"374\u00f77=\", \"291\u00f76=\"),\n    @(\"375\u00f79=\", \"530\u00f72=\"),\n    @(\"902\u00f74=\", \"470\u00f72=\"),\n    @(\"262\u00f77=\", \"783\u00f73=\"),\n    @(\"281\u00f72=\", \"959\u00f76=\"),\n    @(\"413\u00f75=\", \"308\u00f78=\"),\n    @(\"867\u00f72=\", \"649\u00f72=\"),\n    @(\"396\u00f76=\", \"281\u00f78=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\""}
